# Update dashboards - 2026-01-26
# Refresh FRED-sourced series: NGDP level/QoQ SAAR (row 5/6), and the
# Durable Goods Orders (DGORDER) / Non-Defense Capital Goods x Aircraft
# (ADXDNO) M/M & Y/Y blocks (rows 28-31) roll forward one month, from the
# 2025-10-01 release to the 2025-11-01 release.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5 / 6 : NGDP (NGDPSAXDCUSQ) most-recent-quarter figures revised ---
$ws.Range("F5").Value = 7774506.8
$ws.Range("F6").Value = 0.02008474181421249

# --- Row 28 : Dur. Order (DGORDER), M/M % Delta ---
# New release date (2025-11-01) and the value highlighted as newest.
$ws.Range("C28").Value = 45962
$ws.Range("C28").Interior.Color = 65535
$ws.Range("F28").Value = 0.05330084643761013
$ws.Range("G28").Value = -0.02145810719185604
$ws.Range("H28").Value = 0.006436255758670795
$ws.Range("I28").Value = 0.03004963172206243
$ws.Range("J28").Value = -0.02799901206372835

# --- Row 29 : DGORDER, Y/Y % Delta ---
$ws.Range("C29").Value = 45962
$ws.Range("C29").Interior.Color = 65535
$ws.Range("F29").Value = 0.1229486023444545
$ws.Range("G29").Value = 0.04821561312937742
$ws.Range("H29").Value = 0.07412067603746038
$ws.Range("I29").Value = 0.07661265288383932
$ws.Range("J29").Value = 0.03341358778313566

# --- Row 30 : Dur Orders Non Def x Aircraft (ADXDNO), M/M % Delta ---
$ws.Range("C30").Value = 45962
$ws.Range("C30").Interior.Color = 65535
$ws.Range("F30").Value = 0.06566119548130511
$ws.Range("G30").Value = -0.01337070344068647
$ws.Range("H30").Value = 0.001174064535676367
$ws.Range("I30").Value = 0.01907672443132968
$ws.Range("J30").Value = -0.02404555711932721

# --- Row 31 : ADXDNO, Y/Y % Delta ---
$ws.Range("C31").Value = 45962
$ws.Range("C31").Interior.Color = 65535
$ws.Range("F31").Value = 0.1256397039873348
$ws.Range("G31").Value = 0.04954381497984299
$ws.Range("H31").Value = 0.06502168244015354
$ws.Range("I31").Value = 0.06671073894520346
$ws.Range("J31").Value = 0.0329297153895499
